$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 1.82 = 6743.17 pesos
✅ 6743.17 pesos = 1.81 = 822.79 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@
$ws1.Range("A1").Value = $newText

$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 549
$ws2.Range("O10").Value = 3702
$ws2.Range("N12").Value = 3728.94
$ws2.Range("O12").Value = 455
